$wb = $excel.ActiveWorkbook
$wsRaw = $wb.Worksheets.Item("RawData")
$wsBest = $wb.Worksheets.Item("BestPath")
$wsMeta = $wb.Worksheets.Item("Meta")

# =========================================================================
# RawData sheet: recomputed PheromoneScore / Pheromone_Score / Risk_Index /
# Status_Zona / Radius_Visual_KM columns (BI:BM) for rows 2-22
# =========================================================================
$wsRaw.Range("BI2").Value = 0.8975343106629416
$wsRaw.Range("BJ2").Value = 0.8975343106629416
$wsRaw.Range("BK2").Value = 89.75
$wsRaw.Range("BM2").Value = 5.782504626231201

$wsRaw.Range("BI3").Value = 0.9259739153066067
$wsRaw.Range("BJ3").Value = 0.9259739153066067
$wsRaw.Range("BK3").Value = 92.59999999999999
$wsRaw.Range("BM3").Value = 5.839260518582497

$wsRaw.Range("BI4").Value = 0.251511556867749
$wsRaw.Range("BJ4").Value = 0.251511556867749
$wsRaw.Range("BK4").Value = 25.15

$wsRaw.Range("BI5").Value = 0.0001
$wsRaw.Range("BJ5").Value = 0.0001
$wsRaw.Range("BK5").Value = 0.01

$wsRaw.Range("BI6").Value = 0.8548821979343928
$wsRaw.Range("BJ6").Value = 0.8548821979343928
$wsRaw.Range("BK6").Value = 85.48999999999999
$wsRaw.Range("BM6").Value = 3.978392998351124

$wsRaw.Range("BI7").Value = 0.7395082100789241
$wsRaw.Range("BJ7").Value = 0.7395082100789241
$wsRaw.Range("BK7").Value = 73.95
$wsRaw.Range("BM7").Value = 3.817614712715328

$wsRaw.Range("BI8").Value = 0.761243694776387
$wsRaw.Range("BJ8").Value = 0.761243694776387
$wsRaw.Range("BK8").Value = 76.12
$wsRaw.Range("BM8").Value = 3.847903983564655

$wsRaw.Range("BI9").Value = 0.8626015846940868
$wsRaw.Range("BJ9").Value = 0.8626015846940868
$wsRaw.Range("BK9").Value = 86.26000000000001
$wsRaw.Range("BM9").Value = 3.989150273820687

$wsRaw.Range("BI10").Value = 0.8342082463533553
$wsRaw.Range("BJ10").Value = 0.8342082463533553
$wsRaw.Range("BK10").Value = 83.42
$wsRaw.Range("BM10").Value = 5.910978678224625

$wsRaw.Range("BI11").Value = 0.7139596494856451
$wsRaw.Range("BJ11").Value = 0.7139596494856451
$wsRaw.Range("BK11").Value = 71.40000000000001
$wsRaw.Range("BM11").Value = 5.660190157978803

$wsRaw.Range("BI12").Value = 0.9199718003314667
$wsRaw.Range("BJ12").Value = 0.9199718003314667
$wsRaw.Range("BK12").Value = 92
$wsRaw.Range("BM12").Value = 4.069097973368776

$wsRaw.Range("BI13").Value = 0.6800692474526799
$wsRaw.Range("BJ13").Value = 0.6800692474526799
$wsRaw.Range("BK13").Value = 68.01000000000001
$wsRaw.Range("BM13").Value = 3.299552281378881

$wsRaw.Range("BI15").Value = 0.5563274019391721
$wsRaw.Range("BJ15").Value = 0.5563274019391721
$wsRaw.Range("BK15").Value = 55.63
$wsRaw.Range("BM15").Value = 5.101570315656695

$wsRaw.Range("BI16").Value = 0.6068919016003881
$wsRaw.Range("BJ16").Value = 0.6068919016003881
$wsRaw.Range("BK16").Value = 60.69
$wsRaw.Range("BM16").Value = 5.202480062311999

$wsRaw.Range("BI17").Value = 0.7435041105165291
$wsRaw.Range("BJ17").Value = 0.7435041105165291
$wsRaw.Range("BK17").Value = 74.34999999999999
$wsRaw.Range("BL17").Value = "Terdampak"
$wsRaw.Range("BM17").Value = 5.475112116106906

$wsRaw.Range("BI18").Value = 0.5677565934530615
$wsRaw.Range("BJ18").Value = 0.5677565934530615
$wsRaw.Range("BK18").Value = 56.78
$wsRaw.Range("BM18").Value = 5.124379140580679

$wsRaw.Range("BI19").Value = 1
$wsRaw.Range("BJ19").Value = 1
$wsRaw.Range("BK19").Value = 100
$wsRaw.Range("BM19").Value = 4.429323962118182

$wsRaw.Range("BI20").Value = 0.8992978487504288
$wsRaw.Range("BJ20").Value = 0.8992978487504288
$wsRaw.Range("BK20").Value = 89.93000000000001
$wsRaw.Range("BM20").Value = 4.040287991549645

$wsRaw.Range("BI22").Value = 0.9366276693874802
$wsRaw.Range("BJ22").Value = 0.9366276693874802
$wsRaw.Range("BK22").Value = 93.66
$wsRaw.Range("BM22").Value = 4.832155505921507

# =========================================================================
# BestPath sheet: same recomputed columns (BI:BM) for rows 2-22
# =========================================================================
$wsBest.Range("BI2").Value = 0.7435041105165291
$wsBest.Range("BJ2").Value = 0.7435041105165291
$wsBest.Range("BK2").Value = 74.34999999999999
$wsBest.Range("BL2").Value = "Terdampak"
$wsBest.Range("BM2").Value = 5.475112116106906

$wsBest.Range("BI3").Value = 0.7395082100789241
$wsBest.Range("BJ3").Value = 0.7395082100789241
$wsBest.Range("BK3").Value = 73.95
$wsBest.Range("BM3").Value = 3.817614712715328

$wsBest.Range("BI4").Value = 0.251511556867749
$wsBest.Range("BJ4").Value = 0.251511556867749
$wsBest.Range("BK4").Value = 25.15

$wsBest.Range("BI5").Value = 0.8548821979343928
$wsBest.Range("BJ5").Value = 0.8548821979343928
$wsBest.Range("BK5").Value = 85.48999999999999
$wsBest.Range("BM5").Value = 3.978392998351124

$wsBest.Range("BI6").Value = 0.7139596494856451
$wsBest.Range("BJ6").Value = 0.7139596494856451
$wsBest.Range("BK6").Value = 71.40000000000001
$wsBest.Range("BM6").Value = 5.660190157978803

$wsBest.Range("BI8").Value = 0.9366276693874802
$wsBest.Range("BJ8").Value = 0.9366276693874802
$wsBest.Range("BK8").Value = 93.66
$wsBest.Range("BM8").Value = 4.832155505921507

$wsBest.Range("BI9").Value = 0.9259739153066067
$wsBest.Range("BJ9").Value = 0.9259739153066067
$wsBest.Range("BK9").Value = 92.59999999999999
$wsBest.Range("BM9").Value = 5.839260518582497

$wsBest.Range("BI10").Value = 0.0001
$wsBest.Range("BJ10").Value = 0.0001
$wsBest.Range("BK10").Value = 0.01

$wsBest.Range("BI11").Value = 0.8626015846940868
$wsBest.Range("BJ11").Value = 0.8626015846940868
$wsBest.Range("BK11").Value = 86.26000000000001
$wsBest.Range("BM11").Value = 3.989150273820687

$wsBest.Range("BI12").Value = 0.5563274019391721
$wsBest.Range("BJ12").Value = 0.5563274019391721
$wsBest.Range("BK12").Value = 55.63
$wsBest.Range("BM12").Value = 5.101570315656695

$wsBest.Range("BI13").Value = 0.6800692474526799
$wsBest.Range("BJ13").Value = 0.6800692474526799
$wsBest.Range("BK13").Value = 68.01000000000001
$wsBest.Range("BM13").Value = 3.299552281378881

$wsBest.Range("BI14").Value = 0.8342082463533553
$wsBest.Range("BJ14").Value = 0.8342082463533553
$wsBest.Range("BK14").Value = 83.42
$wsBest.Range("BM14").Value = 5.910978678224625

$wsBest.Range("BI15").Value = 0.761243694776387
$wsBest.Range("BJ15").Value = 0.761243694776387
$wsBest.Range("BK15").Value = 76.12
$wsBest.Range("BM15").Value = 3.847903983564655

$wsBest.Range("BI16").Value = 0.9199718003314667
$wsBest.Range("BJ16").Value = 0.9199718003314667
$wsBest.Range("BK16").Value = 92
$wsBest.Range("BM16").Value = 4.069097973368776

$wsBest.Range("BI17").Value = 0.8975343106629416
$wsBest.Range("BJ17").Value = 0.8975343106629416
$wsBest.Range("BK17").Value = 89.75
$wsBest.Range("BM17").Value = 5.782504626231201

$wsBest.Range("BI18").Value = 0.5677565934530615
$wsBest.Range("BJ18").Value = 0.5677565934530615
$wsBest.Range("BK18").Value = 56.78
$wsBest.Range("BM18").Value = 5.124379140580679

$wsBest.Range("BI19").Value = 0.8992978487504288
$wsBest.Range("BJ19").Value = 0.8992978487504288
$wsBest.Range("BK19").Value = 89.93000000000001
$wsBest.Range("BM19").Value = 4.040287991549645

$wsBest.Range("BI21").Value = 1
$wsBest.Range("BJ21").Value = 1
$wsBest.Range("BK21").Value = 100
$wsBest.Range("BM21").Value = 4.429323962118182

$wsBest.Range("BI22").Value = 0.6068919016003881
$wsBest.Range("BJ22").Value = 0.6068919016003881
$wsBest.Range("BK22").Value = 60.69
$wsBest.Range("BM22").Value = 5.202480062311999

# =========================================================================
# Meta sheet: the LSTM/naive prediction columns were trimmed from 10 columns
# (Timestamp, Best_Cost, Node_Count, PredictedLat, PredictedLon,
#  PredictedBearing, PredictedDistanceKM, PredictedDirection,
#  PredictedConfidence, MovementScale) down to 4
# (Timestamp, PredictedBearing, PredictedDistanceKM, PredictedDirection),
# with refreshed values for the single data row.
# =========================================================================
$wsMeta.Range("B1").Value = "PredictedBearing"
$wsMeta.Range("C1").Value = "PredictedDistanceKM"
$wsMeta.Range("D1").Value = "PredictedDirection"

$wsMeta.Range("A2").Value = "2025-12-30T18:12:07.355189"
$wsMeta.Range("B2").Value = 60.72043401770543
$wsMeta.Range("C2").Value = 6.134758280905867
$wsMeta.Range("D2").Value = "NE"

# Drop the now-unused trailing columns (E:J) so the sheet's used range
# shrinks back down to A1:D2, matching the trimmed schema.
$wsMeta.Range("E1:J2").Clear()
